# The commit swaps the contents of ppt/theme/theme1.xml and ppt/theme/theme2.xml:
#   - theme2.xml (the theme actually used by the slide master / all slides,
#     previously the colourful "Integral" theme) becomes the plain default
#     "Office Theme" colour scheme.
#   - theme1.xml (only used by the notes master) would become "Integral",
#     but that part isn't reachable from the slide-facing object model, so
#     we focus on the visible, reachable effect: recolouring the theme that
#     drives the deck's slides via the modern ThemeColorScheme API, which
#     writes straight through to the clrScheme of the shared theme part.

$p = $ppt.ActivePresentation

function New-ComRgb($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target values = the "Office Theme" colour scheme (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) in the same 1-12 index order exposed by
# ThemeColorScheme.Colors(i).
$officeThemeColors = @(
    (New-ComRgb 0x00 0x00 0x00), # 1  dk1      000000
    (New-ComRgb 0xFF 0xFF 0xFF), # 2  lt1      FFFFFF
    (New-ComRgb 0x44 0x54 0x6A), # 3  dk2      44546A
    (New-ComRgb 0xE7 0xE6 0xE6), # 4  lt2      E7E6E6
    (New-ComRgb 0x5B 0x9B 0xD5), # 5  accent1  5B9BD5
    (New-ComRgb 0xED 0x7D 0x31), # 6  accent2  ED7D31
    (New-ComRgb 0xA5 0xA5 0xA5), # 7  accent3  A5A5A5
    (New-ComRgb 0xFF 0xC0 0x00), # 8  accent4  FFC000
    (New-ComRgb 0x44 0x72 0xC4), # 9  accent5  4472C4
    (New-ComRgb 0x70 0xAD 0x47), # 10 accent6  70AD47
    (New-ComRgb 0x05 0x63 0xC1), # 11 hlink    0563C1
    (New-ComRgb 0x95 0x4F 0x72)  # 12 folHlink 954F72
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
